$wb = $excel.ActiveWorkbook

# Rename the "Type" sheet to "ProductType"
$typeSheet = $wb.Worksheets.Item("Type")
$typeSheet.Name = "ProductType"

# The "Company" sheet's view should no longer be the tab-selected one, and
# should be scrolled so row 7 is at the top of the visible window (best
# effort - leaves the existing B2 selection untouched).
$companySheet = $wb.Worksheets.Item("Company")
$companySheet.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# Make "ProductType" the active / tab-selected sheet (activeTab=3).
$typeSheet.Activate()
